$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range('D2').Value = '66.733.07'
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').Value = '3.064.77'
$ws.Range('E3').Value = '  -1.64%  '
$ws.Range('E4').Value = '  +0.04%  '
Set-TextValue ($ws.Range('D5')) '573.76'
$ws.Range('E5').Value = '  -0.58%  '
Set-TextValue ($ws.Range('D6')) '168.32'
$ws.Range('E6').Value = '  -1.85%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '3.062.90'
$ws.Range('E8').Value = '  -1.56%  '
Set-TextValue ($ws.Range('D9')) '0.510'
$ws.Range('E9').Value = '  -2.16%  '
Set-TextValue ($ws.Range('D10')) '6.37'
$ws.Range('E10').Value = '  -2.06%  '
$ws.Range('E11').Value = '  -1.92%  '
Set-TextValue ($ws.Range('D12')) '0.466'
$ws.Range('E12').Value = '  -3.56%  '
$ws.Range('E13').Value = '  -2.75%  '
Set-TextValue ($ws.Range('D14')) '35.50'
$ws.Range('E14').Value = '  -4.47%  '
$ws.Range('E15').Value = '  -1.91%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '66.630.51'
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('B17').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C17').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D17').Value = '3.572.13'
$ws.Range('E17').Value = '  -1.55%  '
Set-TextValue ($ws.Range('D18')) '6.98'
$ws.Range('E18').Value = '  -2.19%  '
Set-TextValue ($ws.Range('D19')) '16.77'
$ws.Range('E19').Value = '  +2.55%  '
$ws.Range('D20').Value = '3.053.72'
$ws.Range('E20').Value = '  -1.91%  '
Set-TextValue ($ws.Range('D21')) '489.97'
$ws.Range('E21').Value = '  +2.79%  '
$ws.Range('E22').Value = '  -3.83%  '
Set-TextValue ($ws.Range('D23')) '7.68'
$ws.Range('E23').Value = '  -3.22%  '
Set-TextValue ($ws.Range('D24')) '82.70'
$ws.Range('E24').Value = '  -1.53%  '
Set-TextValue ($ws.Range('D25')) '12.64'
$ws.Range('E25').Value = '  -6.20%  '
$ws.Range('E26').Value = '  -4.57%  '
$ws.Range('E27').Value = '  +0.70%  '
$ws.Range('E28').Value = '  +0.03%  '
Set-TextValue ($ws.Range('D29')) '7.75'
$ws.Range('E29').Value = '  -2.23%  '
$ws.Range('E30').Value = '  -5.65%  '
$ws.Range('E31').Value = '  -2.58%  '
Set-TextValue ($ws.Range('D32')) '27.43'
$ws.Range('E32').Value = '  -4.01%  '
Set-TextValue ($ws.Range('D33')) '0.111'
$ws.Range('E33').Value = '  -3.61%  '
$ws.Range('D34').Value = '0.0₃0907'
$ws.Range('E34').Value = '  -3.54%  '
Set-TextValue ($ws.Range('D35')) '1.00'
$ws.Range('E35').Value = '  +0.06%  '
Set-TextValue ($ws.Range('D36')) '0.948'
$ws.Range('E36').Value = '  -3.00%  '
Set-TextValue ($ws.Range('D37')) '5.57'
$ws.Range('E37').Value = '  -5.08%  '
Set-TextValue ($ws.Range('D38')) '46.45'
$ws.Range('E38').Value = '  -1.82%  '
$ws.Range('E39').Value = '  -0.20%  '
Set-TextValue ($ws.Range('D40')) '1.96'
$ws.Range('E40').Value = '  -5.47%  '
Set-TextValue ($ws.Range('D41')) '0.300'
$ws.Range('E41').Value = '  -3.63%  '
$ws.Range('E42').Value = '  -4.88%  '
$ws.Range('D43').Value = '2.752.48'
$ws.Range('E43').Value = '  -1.86%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue ($ws.Range('D44')) '135.91'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue ($ws.Range('D45')) '0.0344'
$ws.Range('E45').Value = '  -3.29%  '
Set-TextValue ($ws.Range('D46')) '366.52'
$ws.Range('E46').Value = '  -3.92%  '
$ws.Range('E47').Value = '  -4.89%  '
Set-TextValue ($ws.Range('D49')) '24.36'
$ws.Range('E49').Value = '  -1.91%  '
$ws.Range('E50').Value = '  -2.28%  '
$ws.Range('E51').Value = '  -2.11%  '
